$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Epo"
$ws.Range("C2").Value = "Crlf3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.419201
$ws.Range("H2").Value = 1.257603
$ws.Range("I2").Value = 0.8650221929663464
$ws.Range("J2").Value = 0.8650221929663463
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.866902
$ws.Range("N2").Value = 32.600706
$ws.Range("O2").Value = 0.2279134186299707
$ws.Range("P2").Value = 0.2279134186299707
$ws.Range("Q2").Value = 4.555416185302001
$ws.Range("R2").Value = 40.998745667718
$ws.Range("S2").Value = 0.1971501651897542
$ws.Range("T2").Value = 0.1971501651897541

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Epo"
$ws.Range("C3").Value = "Crlf3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.419201
$ws.Range("H3").Value = 1.257603
$ws.Range("I3").Value = 0.8650221929663464
$ws.Range("J3").Value = 0.8650221929663463
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.066838000000001
$ws.Range("N3").Value = 18.200514
$ws.Range("O3").Value = 0.1272408446173725
$ws.Range("P3").Value = 0.1272408446173725
$ws.Range("Q3").Value = 2.543224556438
$ws.Range("R3").Value = 22.889021007942
$ws.Range("S3").Value = 0.1100661544458097
$ws.Range("T3").Value = 0.1100661544458096

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Epo"
$ws.Range("C4").Value = "Crlf3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.419201
$ws.Range("H4").Value = 1.257603
$ws.Range("I4").Value = 0.8650221929663464
$ws.Range("J4").Value = 0.8650221929663463
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.801689
$ws.Range("N4").Value = 68.405067
$ws.Range("O4").Value = 0.4782237744048302
$ws.Range("P4").Value = 0.4782237744048301
$ws.Range("Q4").Value = 9.558490830489
$ws.Range("R4").Value = 86.026417474401
$ws.Range("S4").Value = 0.4136741780643095
$ws.Range("T4").Value = 0.4136741780643094

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Epo"
$ws.Range("C5").Value = "Crlf3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.419201
$ws.Range("H5").Value = 1.257603
$ws.Range("I5").Value = 0.8650221929663464
$ws.Range("J5").Value = 0.8650221929663463
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.944527999999999
$ws.Range("N5").Value = 23.833584
$ws.Range("O5").Value = 0.1666219623478268
$ws.Range("P5").Value = 0.1666219623478268
$ws.Range("Q5").Value = 3.330354082128
$ws.Range("R5").Value = 29.973186739152
$ws.Range("S5").Value = 0.1441316952664731
$ws.Range("T5").Value = 0.1441316952664731

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Epo"
$ws.Range("C6").Value = "Crlf3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.065412
$ws.Range("H6").Value = 0.196236
$ws.Range("I6").Value = 0.1349778070336537
$ws.Range("J6").Value = 0.1349778070336536
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.866902
$ws.Range("N6").Value = 32.600706
$ws.Range("O6").Value = 0.2279134186299707
$ws.Range("P6").Value = 0.2279134186299707
$ws.Range("Q6").Value = 0.7108257936240001
$ws.Range("R6").Value = 6.397432142616
$ws.Range("S6").Value = 0.03076325344021651
$ws.Range("T6").Value = 0.0307632534402165

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Epo"
$ws.Range("C7").Value = "Crlf3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.065412
$ws.Range("H7").Value = 0.196236
$ws.Range("I7").Value = 0.1349778070336537
$ws.Range("J7").Value = 0.1349778070336536
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.066838000000001
$ws.Range("N7").Value = 18.200514
$ws.Range("O7").Value = 0.1272408446173725
$ws.Range("P7").Value = 0.1272408446173725
$ws.Range("Q7").Value = 0.3968440072560001
$ws.Range("R7").Value = 3.571596065304
$ws.Range("S7").Value = 0.01717469017156281
$ws.Range("T7").Value = 0.01717469017156281

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Epo"
$ws.Range("C8").Value = "Crlf3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.065412
$ws.Range("H8").Value = 0.196236
$ws.Range("I8").Value = 0.1349778070336537
$ws.Range("J8").Value = 0.1349778070336536
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.801689
$ws.Range("N8").Value = 68.405067
$ws.Range("O8").Value = 0.4782237744048302
$ws.Range("P8").Value = 0.4782237744048301
$ws.Range("Q8").Value = 1.491504080868
$ws.Range("R8").Value = 13.423536727812
$ws.Range("S8").Value = 0.0645495963405207
$ws.Range("T8").Value = 0.06454959634052067

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Epo"
$ws.Range("C9").Value = "Crlf3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.065412
$ws.Range("H9").Value = 0.196236
$ws.Range("I9").Value = 0.1349778070336537
$ws.Range("J9").Value = 0.1349778070336536
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.944527999999999
$ws.Range("N9").Value = 23.833584
$ws.Range("O9").Value = 0.1666219623478268
$ws.Range("P9").Value = 0.1666219623478268
$ws.Range("Q9").Value = 0.519667465536
$ws.Range("R9").Value = 4.677007189824
$ws.Range("S9").Value = 0.02249026708135367
$ws.Range("T9").Value = 0.02249026708135366

